$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column before J ("Additional Comments" shifts from J to K)
$ws.Columns.Item(10).Insert()

# Move the existing "Additional Comments" comment from J1 (post-shift) to K1,
# and put the new "Due Date" header + comment in J1.
$oldComment = $ws.Range("J1").Comment
$oldText = $oldComment.Text()
$oldComment.Delete()

$ws.Range("J1").Value = "Due Date"
$ws.Range("K1").AddComment($oldText)
$ws.Range("J1").AddComment("If this ticket has a due date, enter it here in YYYY-MM-DD format.")

# Touch L1 so the sheet's used range / dimension extends to include it.
$ws.Range("L1").NumberFormat = "General"

# Match the resulting selection from the authored edit.
$ws.Activate()
$ws.Range("K9").Select()
